$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 303
$ws.Range("B2").Value = 6
$ws.Range("C2").Value = "王*佑"
$ws.Range("D2").Value = "2023-07-05 15:19:43"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "IN"
